$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the data range so numeric-looking strings
# (e.g. "64.318.06", "1.00") are not auto-converted to numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "64.318.06"
$ws.Range("D3").Value = "3.503.61"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "586.79"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "134.55"
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "4.100.95"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").Value = "3.504.39"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "26.03"
$ws.Range("E16").Value = "  -4.69%  "
$ws.Range("D17").Value = "64.330.63"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "5.75"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "13.74"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").Value = "395.01"
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").Value = "0.573"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "3.643.79"
$ws.Range("D24").Value = "74.29"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").Value = "0.0000116"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "7.39"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").Value = "1.50"
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "3.522.40"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "0.150"
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "23.48"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "161.61"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").Value = "0.807"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "25.48"
$ws.Range("E43").Value = "  -4.79%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "4.43"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "1.65"
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").Value = "1.17"
$ws.Range("E47").Value = "  -2.93%  "
$ws.Range("D48").Value = "2.468.04"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "0.896"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  -1.17%  "

# Restore original (default) cell formatting now that text values are set,
# so style indices match the unedited cells.
$ws.Range("D2:E51").ClearFormats()
